$wb = $excel.ActiveWorkbook

# --- Sheet 1: Stocks ---
$ws1 = $wb.Worksheets.Item("Stocks")
$ws1.Cells.Item(2, 2).Value = 0.0006421214165380263
$ws1.Cells.Item(2, 3).Value = 0.022163179588173
$ws1.Cells.Item(3, 2).Value = 0.0009525673071090029
$ws1.Cells.Item(3, 3).Value = 0.01763312565338153
$ws1.Cells.Item(4, 2).Value = 0.0005293197922399937
$ws1.Cells.Item(4, 3).Value = 0.02213256403877736
$ws1.Cells.Item(5, 2).Value = 0.0003449730979418267
$ws1.Cells.Item(5, 3).Value = 0.02162974797326941
$ws1.Cells.Item(6, 2).Value = 0.001339861789141646
$ws1.Cells.Item(6, 3).Value = 0.02401398385755637
$ws1.Cells.Item(9, 2).Value = 0.0008594941721883319
$ws1.Cells.Item(9, 3).Value = 0.01909794576836636
$ws1.Cells.Item(10, 2).Value = 0.001202652772747439
$ws1.Cells.Item(10, 3).Value = 0.01827319193348403
$ws1.Cells.Item(11, 2).Value = -0.0004646364889073274
$ws1.Cells.Item(11, 3).Value = 0.02150403437377608

# --- Sheet 2: Crypto ---
$ws2 = $wb.Worksheets.Item("Crypto")
$ws2.Cells.Item(2, 2).Value = -0.0009462549077120909
$ws2.Cells.Item(2, 3).Value = 0.03630537429437645
$ws2.Cells.Item(3, 2).Value = 0.0004004341971715087
$ws2.Cells.Item(3, 3).Value = 0.04330960451409918
$ws2.Cells.Item(4, 2).Value = -0.0009475469251185598
$ws2.Cells.Item(4, 3).Value = 0.04485220147448776
$ws2.Cells.Item(5, 2).Value = 0.001407352743383836
$ws2.Cells.Item(5, 3).Value = 0.04430072809446627
$ws2.Cells.Item(6, 2).Value = -0.0002666370955085417
$ws2.Cells.Item(6, 3).Value = 0.02895642791321457
$ws2.Cells.Item(7, 2).Value = 0.0001698017340538461
$ws2.Cells.Item(7, 3).Value = 0.02728721768065267
$ws2.Cells.Item(8, 2).Value = -0.001935140215022296
$ws2.Cells.Item(8, 3).Value = 0.03557390366834624
$ws2.Cells.Item(9, 2).Value = 0.0005055873766004
$ws2.Cells.Item(9, 3).Value = 0.04694877812707692
$ws2.Cells.Item(10, 2).Value = -0.001076934762782215
$ws2.Cells.Item(10, 3).Value = 0.03583406844817972
$ws2.Cells.Item(11, 2).Value = 0.0005457174036354702
$ws2.Cells.Item(11, 3).Value = 0.04918708106781185
$ws2.Cells.Item(12, 2).Value = 0.0004779069833942835
$ws2.Cells.Item(12, 3).Value = 0.03619237636843207
$ws2.Cells.Item(13, 1).Value = "KAS-USD"
$ws2.Cells.Item(13, 2).Value = 0.01744103067694832
$ws2.Cells.Item(13, 3).Value = 0.1230996776299637
$ws2.Cells.Item(14, 2).Value = -0.0002853039963843259
$ws2.Cells.Item(14, 3).Value = 0.02466974165848193
$ws2.Cells.Item(15, 2).Value = 0.001232268108914391
$ws2.Cells.Item(15, 3).Value = 0.04220231508622439
$ws2.Cells.Item(16, 2).Value = 0.0009022075568424984
$ws2.Cells.Item(16, 3).Value = 0.04038740999764472
$ws2.Cells.Item(17, 2).Value = 0.0009886723621352019
$ws2.Cells.Item(17, 3).Value = 0.05065994433348981
$ws2.Cells.Item(18, 2).Value = 0.003389285708800178
$ws2.Cells.Item(18, 3).Value = 0.03770878311670382
$ws2.Cells.Item(19, 2).Value = 0.0004775359105210174
$ws2.Cells.Item(19, 3).Value = 0.05382868312718332
$ws2.Cells.Item(20, 2).Value = 0.0002351162794294515
$ws2.Cells.Item(20, 3).Value = 0.05479576322261236
$ws2.Cells.Item(21, 2).Value = -0.0001231138256954889
$ws2.Cells.Item(21, 3).Value = 0.04772686211557989
$ws2.Cells.Item(22, 2).Value = 0.0004388946013284949
$ws2.Cells.Item(22, 3).Value = 0.02490479007018209
$ws2.Cells.Item(23, 2).Value = 0.000179030553937395
$ws2.Cells.Item(23, 3).Value = 0.02758758242718384
$ws2.Cells.Item(24, 2).Value = 0.0004373874387499149
$ws2.Cells.Item(24, 3).Value = 0.04372835900203738
$ws2.Cells.Item(25, 2).Value = -0.0001819086296563999
$ws2.Cells.Item(25, 3).Value = 0.02952524701840095
$ws2.Cells.Item(26, 2).Value = 0.001573038748370876
$ws2.Cells.Item(26, 3).Value = 0.04954929129429662

Write-Host "Edit applied successfully"
